# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    carries the localization status column.
# 2) Narrow the "zh-cn"/"de-de" status columns (Overview!E:F, zh-cn!C,
#    de-de!C) from ~17.22 chars down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Update status text wherever it appears ---------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2) Narrow the status columns -----------------------------------------
# Stored column width 17.2159881591797 chars -> 13.4101845877511 chars.
# The ColumnWidth property (in characters) is what drives the stored
# <col width="..."/> attribute, so set it directly on each column.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
